# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 5, pushing the existing rows
# 5-48 down to 6-49 (dates/values shift down by one row accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 5 - shifts rows 5..48 down to 6..49
# and inherits the formatting (e.g. the date-formatted style) of the row
# that used to be at position 5.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44685
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 100112022
$ws.Cells.Item(5, 7).Value = "Arveja Verde"
$ws.Cells.Item(5, 8).Value = "Perfection"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 160
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 27000
$ws.Cells.Item(5, 13).Value = 26000
$ws.Cells.Item(5, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 1040
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
